$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Insert a new row at 21, shifting old rows 21-29 down to 22-30
$ws.Rows.Item(21).Insert()

# Fill in the new row 21 with the new task
$ws.Range("B21").Value = "Skapa metod för introduktionstext"
$ws.Range("C21").Value = "Färdig"
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 2

# Remove E22/E23 (previously E21/E22, shifted down), add F22/F23 instead
$ws.Range("E22").Clear()
$ws.Range("E23").Clear()
$ws.Range("F22").HorizontalAlignment = -4131
$ws.Range("F23").HorizontalAlignment = -4131

# Clear leftover C24 cell from old row23 shift
$ws.Range("C24").Clear()

# Clear leftover row 30 (old row 29 pushed down by insert); table stays at 29 rows
$ws.Range("D30").Clear()

# Update the selection to match target
[void]$ws.Range("D22:D23").Select()
